$wb = $excel.ActiveWorkbook

# "Generate Report for Handback" refreshes the handoff/handback timestamps
# recorded on the Overview sheet and the per-locale (zh-cn / de-de) sheets.

$overview = $wb.Sheets.Item("Overview")
$zhcn = $wb.Sheets.Item("zh-cn")
$dede = $wb.Sheets.Item("de-de")

# Overview!G2 = "Latest HO Xliff Generate Date" for f12df289...md
# de-de!H2   = "Correspond Handoff Datetime" for f12df289...md (de-de)
# These two cells share the same underlying timestamp value.
$overview.Range("G2").Value = "2016-09-02 15:21:44"
$dede.Range("H2").Value = "2016-09-02 15:21:44"

# zh-cn!H2 = "Correspond Handoff Datetime" for f12df289...md (zh-cn)
$zhcn.Range("H2").Value = "2016-09-02 15:21:39"

# zh-cn!K2 = "Correspond Handback DateTime" for f12df289...md (zh-cn)
$zhcn.Range("K2").Value = "2016-09-02 15:21:57"

# de-de!K2 = "Correspond Handback DateTime" for f12df289...md (de-de)
$dede.Range("K2").Value = "2016-09-02 15:22:16"
